$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: perform ALL structural changes first (rename / add / move).
# Worksheet object handles can go stale across Add()/Move() calls in this
# runtime, so we deliberately re-resolve every sheet by name once all the
# structural operations below are finished, before touching any cell data.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Name = "Estoque"

$produtos = $wb.Worksheets.Add()
$produtos.Name = "Produtos"

$vendas = $wb.Worksheets.Add()
$vendas.Name = "Vendas"

$vendas = $wb.Worksheets.Item("Vendas")
$estoque = $wb.Worksheets.Item("Estoque")
$vendas.Move(, $estoque)

# ---------------------------------------------------------------------------
# Step 2: re-fetch fresh references now that the sheet collection is stable:
# Produtos, Estoque, Vendas (in this tab order).
# ---------------------------------------------------------------------------
$produtos = $wb.Worksheets.Item("Produtos")
$estoque = $wb.Worksheets.Item("Estoque")
$vendas = $wb.Worksheets.Item("Vendas")

# ---------------------------------------------------------------------------
# Helper formatting: reproduce the bold / centered / bordered header style
# ---------------------------------------------------------------------------
function Set-HeaderFormat($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# Step 3: Produtos sheet data
# ---------------------------------------------------------------------------
$produtosHeader = $produtos.Range("A1:E1")
Set-HeaderFormat($produtosHeader)
$produtos.Range("A1").Value = "Produto"
$produtos.Range("B1").Value = "Marca"
$produtos.Range("C1").Value = "Método_Compra"
$produtos.Range("D1").Value = "Valor_Método"
$produtos.Range("E1").Value = "Método_Venda"

$produtos.Range("A2").Value = "Ração"
$produtos.Range("B2").Value = "Pedigree"
$produtos.Range("C2").Value = "Pacote"
$produtos.Range("D2").Value = 20

$produtos.Range("A3").Value = "Shampoo"
$produtos.Range("B3").Value = "Gato"
$produtos.Range("C3").Value = "Unidade"
$produtos.Range("D3").NumberFormat = "@"
$produtos.Range("D3").Value = "2"
$produtos.Range("E3").Value = "Unidade"

# ---------------------------------------------------------------------------
# Step 4: Estoque sheet data
# ---------------------------------------------------------------------------
$estoque.Range("A1").Value = "Produto"
$estoque.Range("B1").Value = "Marca"
$estoque.Range("C1").Value = "Quantidade"
$estoque.Range("D1").Value = "Valor Total Gasto"

$estoque.Range("A2").Value = "Ração"
$estoque.Range("B2:D2").ClearContents()

$estoque.Range("A3").Value = "Ração"
$estoque.Range("B3").Value = "Pedigree"
$estoque.Range("C3").NumberFormat = "@"
$estoque.Range("C3").Value = "6"
$estoque.Range("D3").NumberFormat = "@"
$estoque.Range("D3").Value = "20.0"

$estoque.Range("A4").Value = "Shampoo"
$estoque.Range("B4").Value = "Gato"
$estoque.Range("C4").NumberFormat = "@"
$estoque.Range("C4").Value = "100"
$estoque.Range("D4").NumberFormat = "@"
$estoque.Range("D4").Value = "2.0"

# ---------------------------------------------------------------------------
# Step 5: Vendas sheet data
# ---------------------------------------------------------------------------
$vendasHeader = $vendas.Range("A1:B1")
Set-HeaderFormat($vendasHeader)
$vendas.Range("A1").Value = "Produto"
$vendas.Range("B1").Value = "Marca"

$vendas.Range("A2").Value = "Ração"

# ---------------------------------------------------------------------------
# Step 6: leave Produtos as the active/selected sheet, matching the original
# workbook where the single sheet had tabSelected="1"
# ---------------------------------------------------------------------------
$produtos.Activate()

Write-Host "Sheet order:"
foreach ($s in $wb.Worksheets) { Write-Host " - $($s.Name)" }
